$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.447.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.597.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.33%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.385'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.060.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.239.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("E16").Value = '  +3.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.590.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.15%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("E26").Value = '  -3.79%  '
$ws.Range("E27").Value = '  +1.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '548.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0849'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.80'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("E40").Value = '  -4.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '164.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0577'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("E46").Value = '  +3.79%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0956'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0227'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.68%  '
